$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 7
$ws.Range("A2").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A7").Value = 44449
$ws.Range("E7").Value = "['coz', 'qxb', 'ayt']"
$ws.Range("B7").Value = "model_floodwater_unet_pc_augm_diceloss_2"
$ws.Range("C7").Value = 0.715

# New row 8
$ws.Range("A2").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = 44450
$ws.Range("B8").Value = "model_floodwater_unet_pc_augm_IOUloss"
$ws.Range("C8").Value = 0.663
$ws.Range("E8").Value = "['kuo', 'wvy', 'awc']"

$ws.Range("E6").Select()
